$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 50,4
$arr[0,0] = 0.04872769489884377
$arr[0,1] = 0.9913510084152222
$arr[0,2] = 0.1162698268890381
$arr[0,3] = 0.963495135307312
$arr[1,0] = 0.00778226787224412
$arr[1,1] = 0.998447597026825
$arr[1,2] = 0.06846234947443008
$arr[1,3] = 0.9778631925582886
$arr[2,0] = 0.00463435472920537
$arr[2,1] = 0.9987096786499023
$arr[2,2] = 0.01203029230237007
$arr[2,3] = 0.9985799193382263
$arr[3,0] = 0.001854672096669674
$arr[3,1] = 0.9994354844093323
$arr[3,2] = 0.01542354468256235
$arr[3,3] = 0.9951549768447876
$arr[4,0] = 0.001591858570463955
$arr[4,1] = 0.9996572732925415
$arr[4,2] = 0.00940130278468132
$arr[4,3] = 0.9990811347961426
$arr[5,0] = 0.001263531274162233
$arr[5,1] = 0.9996774196624756
$arr[5,2] = 0.0175401009619236
$arr[5,3] = 0.9929830431938171
$arr[6,0] = 0.0009633513982407749
$arr[6,1] = 0.9996371269226074
$arr[6,2] = 0.001459159655496478
$arr[6,3] = 0.9998329281806946
$arr[7,0] = 0.001066343393176794
$arr[7,1] = 0.9997984170913696
$arr[7,2] = 0.004663672298192978
$arr[7,3] = 0.9991646409034729
$arr[8,0] = 0.0005156368133611977
$arr[8,1] = 0.9998992085456848
$arr[8,2] = 0.002270422410219908
$arr[8,3] = 0.999248206615448
$arr[9,0] = 0.001056618755683303
$arr[9,1] = 0.9997782111167908
$arr[9,2] = 0.001016339985653758
$arr[9,3] = 0.9997494220733643
$arr[10,0] = 0.0003958890156354755
$arr[10,1] = 0.9998588562011719
$arr[10,2] = 0.0008839988731779158
$arr[10,3] = 0.9998329281806946
$arr[11,0] = 0.0005234397249296308
$arr[11,1] = 0.9997782111167908
$arr[11,2] = 0.0009705186239443719
$arr[11,3] = 0.9997494220733643
$arr[12,0] = 0.0004141610697843134
$arr[12,1] = 0.9998992085456848
$arr[12,2] = 0.001528264256194234
$arr[12,3] = 0.9998329281806946
$arr[13,0] = 0.00008285167859867215
$arr[13,1] = 0.9999798536300659
$arr[13,2] = 0.001273514702916145
$arr[13,3] = 0.9998329281806946
$arr[14,0] = 0.0003121593326795846
$arr[14,1] = 0.9999193549156189
$arr[14,2] = 0.001323317643254995
$arr[14,3] = 0.9998329281806946
$arr[15,0] = 0.0001992267643800005
$arr[15,1] = 0.9999193549156189
$arr[15,2] = 0.00167994829826057
$arr[15,3] = 0.9998329281806946
$arr[16,0] = 0.0004221591516397893
$arr[16,1] = 0.9998992085456848
$arr[16,2] = 0.001895668567158282
$arr[16,3] = 0.9996658563613892
$arr[17,0] = 0.0002443742996547371
$arr[17,1] = 0.9999193549156189
$arr[17,2] = 0.002119817305356264
$arr[17,3] = 0.9996658563613892
$arr[18,0] = 0.000206456141313538
$arr[18,1] = 0.9999798536300659
$arr[18,2] = 0.001915992004796863
$arr[18,3] = 0.9998329281806946
$arr[19,0] = 0.0007310992805287242
$arr[19,1] = 0.9998387098312378
$arr[19,2] = 0.002349842339754105
$arr[19,3] = 0.9998329281806946
$arr[20,0] = 0.0003132978163193911
$arr[20,1] = 0.999939501285553
$arr[20,2] = 0.00182623160071671
$arr[20,3] = 0.9996658563613892
$arr[21,0] = 0.00003318466042401269
$arr[21,1] = 0.9999798536300659
$arr[21,2] = 0.00202444102615118
$arr[21,3] = 0.9996658563613892
$arr[22,0] = 0.0009933458641171455
$arr[22,1] = 0.9998185634613037
$arr[22,2] = 0.001152956043370068
$arr[22,3] = 0.9996658563613892
$arr[23,0] = 0.0002190259110648185
$arr[23,1] = 0.999939501285553
$arr[23,2] = 0.002094593364745378
$arr[23,3] = 0.9998329281806946
$arr[24,0] = 0.0005047370214015245
$arr[24,1] = 0.9998790621757507
$arr[24,2] = 0.00181324640288949
$arr[24,3] = 0.9995823502540588
$arr[25,0] = 0.00005787853297078982
$arr[25,1] = 0.9999798536300659
$arr[25,2] = 0.00162134354468435
$arr[25,3] = 0.9995823502540588
$arr[26,0] = 0.00005760549174738117
$arr[26,1] = 0.9999798536300659
$arr[26,2] = 0.003307906445115805
$arr[26,3] = 0.9995823502540588
$arr[27,0] = 0.0003011898952536285
$arr[27,1] = 0.9998992085456848
$arr[27,2] = 0.01857579313218594
$arr[27,3] = 0.9941525459289551
$arr[28,0] = 0.0002523934235796332
$arr[28,1] = 0.9999597072601318
$arr[28,2] = 0.003325604135170579
$arr[28,3] = 0.9996658563613892
$arr[29,0] = 0.000008441307727480307
$arr[29,1] = 1.0
$arr[29,2] = 0.003541434183716774
$arr[29,3] = 0.9996658563613892
$arr[30,0] = 0.0001541435776744038
$arr[30,1] = 0.999939501285553
$arr[30,2] = 0.002098003169521689
$arr[30,3] = 0.9998329281806946
$arr[31,0] = 0.0007424585637636483
$arr[31,1] = 0.9998387098312378
$arr[31,2] = 0.001337770256213844
$arr[31,3] = 0.9998329281806946
$arr[32,0] = 0.000137047391035594
$arr[32,1] = 0.9999597072601318
$arr[32,2] = 0.002440424636006355
$arr[32,3] = 0.9996658563613892
$arr[33,0] = 0.00002120395220117643
$arr[33,1] = 0.9999798536300659
$arr[33,2] = 0.002577274572104216
$arr[33,3] = 0.9996658563613892
$arr[34,0] = 0.0001176949808723293
$arr[34,1] = 0.9999798536300659
$arr[34,2] = 0.002308253897354007
$arr[34,3] = 0.9998329281806946
$arr[35,0] = 0.0004537097120191902
$arr[35,1] = 0.9998588562011719
$arr[35,2] = 0.002078523859381676
$arr[35,3] = 0.9998329281806946
$arr[36,0] = 0.0002284746151417494
$arr[36,1] = 0.999939501285553
$arr[36,2] = 0.00232000183314085
$arr[36,3] = 0.9998329281806946
$arr[37,0] = 0.0002947113825939596
$arr[37,1] = 0.9999193549156189
$arr[37,2] = 0.003107803408056498
$arr[37,3] = 0.9995823502540588
$arr[38,0] = 0.00001482611241954146
$arr[38,1] = 1.0
$arr[38,2] = 0.003810758702456951
$arr[38,3] = 0.9996658563613892
$arr[39,0] = 0.000194270905922167
$arr[39,1] = 0.9999597072601318
$arr[39,2] = 0.003062099916860461
$arr[39,3] = 0.9995823502540588
$arr[40,0] = 0.0001900628121802583
$arr[40,1] = 0.9999193549156189
$arr[40,2] = 0.003283705795183778
$arr[40,3] = 0.9998329281806946
$arr[41,0] = 0.0001322337338933721
$arr[41,1] = 0.9999597072601318
$arr[41,2] = 0.003577649127691984
$arr[41,3] = 0.9998329281806946
$arr[42,0] = 0.0001524569961475208
$arr[42,1] = 0.9999597072601318
$arr[42,2] = 0.001904546981677413
$arr[42,3] = 0.9998329281806946
$arr[43,0] = 0.0002568055060692132
$arr[43,1] = 0.9999597072601318
$arr[43,2] = 0.003277556970715523
$arr[43,3] = 0.9998329281806946
$arr[44,0] = 0.0000113979740490322
$arr[44,1] = 1.0
$arr[44,2] = 0.003416521940380335
$arr[44,3] = 0.9998329281806946
$arr[45,0] = 0.00007801568426657468
$arr[45,1] = 0.9999798536300659
$arr[45,2] = 0.003669193014502525
$arr[45,3] = 0.9998329281806946
$arr[46,0] = 0.000009052690074895509
$arr[46,1] = 1.0
$arr[46,2] = 0.004442702047526836
$arr[46,3] = 0.9998329281806946
$arr[47,0] = 0.0004959235666319728
$arr[47,1] = 0.9998992085456848
$arr[47,2] = 0.004369628615677357
$arr[47,3] = 0.9998329281806946
$arr[48,0] = 0.00007549977453891188
$arr[48,1] = 0.9999597072601318
$arr[48,2] = 0.005136909428983927
$arr[48,3] = 0.9998329281806946
$arr[49,0] = 0.000544364913366735
$arr[49,1] = 0.999939501285553
$arr[49,2] = 0.00380149926058948
$arr[49,3] = 0.9998329281806946
$ws.Range("A2:D51").Value = $arr
